# Lecture 7 notes insertion
# Locate the paragraph "All we care about is the .so files at the end."
# (the last bulleted item before the trailing bookmark paragraph) and
# insert the new Lecture 7 content immediately after it.

$d = $word.ActiveDocument

$anchorText = "All we care about is the .so files at the end."

$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like ($anchorText + "*")) {
        $anchorIndex = $i
    }
}

if ($anchorIndex -eq -1) {
    throw "Could not find anchor paragraph: $anchorText"
}

$anchorPara = $d.Paragraphs.Item($anchorIndex)
$insertAt = $anchorPara.Range.End
$insertRange = $d.Range($insertAt, $insertAt)

# The block below is the OOXML for the five new paragraphs described in the
# commit diff:
#   1. a blank paragraph
#   2. "Lecture 7 - 5/21/13"
#   3. "Threading and Multiprocessing" (underlined heading)
#   4. a bulleted item: "Went over argparse"
#   5. a bulleted item about multithreading/multiprocessing/the GIL
#
# A trailing empty <w:p/> is appended because Word's InsertXML merges the
# paragraph mark of the final inserted paragraph into whatever paragraph
# follows the insertion point; the extra paragraph mark absorbs that merge
# and is removed afterwards, leaving the five target paragraphs intact and
# the original following paragraph (the _GoBack bookmark paragraph)
# untouched.
$newParagraphsXml = @'
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:rPr><w:rFonts w:cs="Courier New"/></w:rPr></w:pPr></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:rPr><w:rFonts w:cs="Courier New"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="Courier New"/></w:rPr><w:t>Lecture 7 – 5/21/13</w:t></w:r></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:rPr><w:rFonts w:cs="Courier New"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="Courier New"/><w:u w:val="single"/></w:rPr><w:t>Threading and Multiprocessing</w:t></w:r></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:cs="Courier New"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="Courier New"/></w:rPr><w:t xml:space="preserve">Went over </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cs="Courier New"/></w:rPr><w:t>argparse</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:cs="Courier New"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="Courier New"/></w:rPr><w:t>Multithreading is good for IO things</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Courier New"/></w:rPr><w:t xml:space="preserve"> –</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Courier New"/></w:rPr><w:t xml:space="preserve"> multiprocessing is necessary otherwise.</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Courier New"/></w:rPr><w:t xml:space="preserve"> The GIL gets in the way.</w:t></w:r></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'></w:p>
'@

$null = $insertRange.InsertXML($newParagraphsXml)

# InsertXML placed 6 new paragraph marks after the anchor paragraph (the 5
# intended ones plus the trailing spacer that absorbed the merge with the
# paragraph that used to directly follow the anchor). Remove that spacer —
# it is the 6th paragraph after the anchor. A paragraph holding no actual
# text still reports its own paragraph-mark character, so compare against
# that rather than an empty string.
$strayIndex = $anchorIndex + 6
$strayPara = $d.Paragraphs.Item($strayIndex)
if ($strayPara.Range.Text.Trim() -eq "") {
    $null = $strayPara.Range.Delete()
} else {
    throw "Unexpected content in paragraph meant to be the InsertXML spacer"
}
